# Generate Report for Handoff
# Refresh the handoff report: rows 4-7 ("Ready for handoff" items) get
# bumped from Priority "low" -> "ht" and their "Latest Handoff Datetime"
# refreshed to the new generation timestamps, on both the zh-cn and
# de-de localization-status sheets.

$wb = $excel.ActiveWorkbook

$zhcn = $wb.Worksheets.Item("zh-cn")
foreach ($r in 4..7) {
    $zhcn.Range("E$r").Value = "ht"
    $zhcn.Range("H$r").Value = "2016-09-02 12:36:35"
}

$dede = $wb.Worksheets.Item("de-de")
foreach ($r in 4..7) {
    $dede.Range("E$r").Value = "ht"
    $dede.Range("H$r").Value = "2016-09-02 12:36:39"
}
